# Update countries & provincias Spain
# - Insert "Trinidad yTobago" ahead of "Guadalupe"/"Aruba" in the country list
#   (rows 136-138 shift: Trinidad yTobago, Guadalupe, Aruba), with refreshed
#   daily statistics for the affected rows.
# - Refresh the "last updated" timestamp string.
# - Refresh numeric COVID statistics (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Muertes hoy, Muertes) for a number of countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) -----------------------------
$ws.Range("A1").Value = "Datos actualizados a 14 de Septiembre de 2020 a las 17:40"

# --- Re-sequence the three rows around the newly inserted country --------
# Row 135 (Sri Lanka) keeps its place; a new "Trinidad yTobago" entry takes
# row 136, pushing the previous Guadalupe/Aruba rows down by one.
$ws.Range("A136").Value = "Trinidad yTobago"
$ws.Range("A137").Value = "Guadalupe"
$ws.Range("A138").Value = "Aruba"

# --- Refresh numeric statistics -------------------------------------------
$updates = @{
    "B4"  = 6716377; "C4"  = 7919;  "D4"  = 3981334; "E4"  = 2536423; "G4"  = 100;  "H4"  = 198620
    "B5"  = 4878042; "C5"  = 33039; "D5"  = 3809549; "E5"  = 988467;  "G5"  = 272;  "H5"  = 80026

    "D13" = 428953;  "E13" = 115172; "G13" = 60;      "H13" = 11412

    "B14" = 436433;  "C14" = 1685;  "D14" = 407725;  "E14" = 16695;  "G14" = 64;   "H14" = 12013

    "B17" = 371125;  "C17" = 2621;  "G17" = 9;        "H17" = 41637

    "B23" = 288761;  "C23" = 1008;  "D23" = 213950;  "E23" = 39187;  "G23" = 14;   "H23" = 35624

    "B29" = 136972;  "C29" = 313;   "D29" = 120564;  "E29" = 7236;   "G29" = 1;    "H29" = 9172

    "B31" = 121975;  "C31" = 235;   "D31" = 118931;  "E31" = 2837;   "G31" = 2;    "H31" = 207

    "B35" = 104110;  "C35" = 450;   "D35" = 77790;   "E35" = 24336;  "G35" = 16;   "H35" = 1984

    "B46" = 80266;   "C46" = 777;   "D46" = 69981;   "E46" = 9886

    "D55" = 56802;   "E55" = 625

    "D61" = 39600;   "E61" = 5811

    "B65" = 43207;   "C65" = 229;   "E65" = 11641;   "G65" = 6;      "H65" = 1129

    "B93" = 12219;   "C93" = 65;    "E93" = 1583

    "B112" = 5396;   "C112" = 1;    "D112" = 5331

    "B135" = 3253;   "C135" = 19;   "E135" = 235

    "B136" = 3091;   "C136" = 49;   "D136" = 787;    "E136" = 2250;  "G136" = 1;   "H136" = 54

    "B137" = 3080;   "D137" = 837;  "E137" = 2219;   "H137" = 24

    "B138" = 3046;   "D138" = 1542; "E138" = 1486;   "H138" = 18

    "B140" = 2974;   "C140" = 46;   "D140" = 1344;   "E140" = 1563

    "B142" = 2872;   "C142" = 67;   "E142" = 1544

    "B144" = 2587;   "C144" = 9;    "E144" = 1248
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
